# ---------------------------------------------------------------------------
# BUR Testing_IPS.xlsx edit script
#
# Summary of changes implemented (per commit message / xml diff):
#  1. Correct the "cvd" label text: "Voluntary Turnover Professional" ->
#     "Professional Voluntary Turnover" (shared across all division sheets).
#  2. Re-run numbers with the frozen 2024 file: several Voluntary-Turnover /
#     Internal-Fill-Rate percentage cells change slightly on the 5 division
#     tabs (Clutches & Brakes, Couplings, Gearing, Industrial Components,
#     Segment Functions).
#  3. The former "L1_IPS" tab is repurposed into a (now mostly blank)
#     "Integration" tab, and a brand new "L1_IPS" tab is appended after it
#     holding the (slightly revised) L1 IPS numbers that used to live on the
#     old tab.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1 - fix the mislabeled "cvd" shared text across every worksheet
# ---------------------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Voluntary Turnover Professional", "Professional Voluntary Turnover")
}

# ---------------------------------------------------------------------------
# Step 2 - split "L1_IPS" into "Integration" (old tab, cleared out) and a new
# "L1_IPS" tab (copy of the old tab, placed right after it) holding the
# refreshed numbers.
# ---------------------------------------------------------------------------
$src = $wb.Worksheets.Item("L1_IPS")
$src.Copy($null, $src)
$newSheet = $wb.Worksheets.Item("L1_IPS (2)")
$src.Name = "Integration"
$newSheet.Name = "L1_IPS"

# --- 2a. Clear out the "Integration" tab's data rows -----------------------
# Rows 2-3 (Professional Voluntary Turnover) keep their label columns
# (A/B/C/E) but the numeric columns are blanked out, except S:V which become
# explicit zeros. Rows 4-7 (the Internal Fill Rate block, and the 3rd
# Professional Voluntary Turnover row) are removed entirely.
$ig = $wb.Worksheets.Item("Integration")
$ig.Range("D2").ClearContents()
$ig.Range("F2:R2").ClearContents()
$ig.Range("D3").ClearContents()
$ig.Range("F3:R3").ClearContents()
$ig.Range("S2:V2").Value = 0
$ig.Range("S3:V3").Value = 0
$ig.Range("A4:V7").Delete()

# --- 2b. Apply the refreshed numbers to the new "L1_IPS" tab ---------------
$l1 = $wb.Worksheets.Item("L1_IPS")

$l1.Range("F2").Value = 0.0113
$l1.Range("G2").Value = 0.0079
$l1.Range("H2").Value = 0.0105
$l1.Range("I2").Value = 0.0297
$l1.Range("J2").Value = 0.0086
$l1.Range("K2").Value = 0.0088
$l1.Range("L2").Value = 0.0104
$l1.Range("M2").Value = 0.0278
$l1.Range("N2").Value = 0.0088
$l1.Range("O2").Value = 0.0098
$l1.Range("P2").Value = 0.007
$l1.Range("Q2").Value = 0.0255
$l1.Range("R2").Value = 0.0068
$l1.Range("S2").Value = 0.0066
$l1.Range("T2").Value = 0.0106
$l1.Range("U2").Value = 0.024
$l1.Range("V2").Value = 0.107

$l1.Range("F3").Value = 0.01017
$l1.Range("G3").Value = 0.00711
$l1.Range("H3").Value = 0.00945
$l1.Range("I3").Value = 0.02673
$l1.Range("J3").Value = 0.00774
$l1.Range("K3").Value = 0.00792
$l1.Range("L3").Value = 0.00936
$l1.Range("M3").Value = 0.02502
$l1.Range("N3").Value = 0.00792
$l1.Range("O3").Value = 0.00882
$l1.Range("P3").Value = 0.0063
$l1.Range("Q3").Value = 0.02295
$l1.Range("R3").Value = 0.00612
$l1.Range("S3").Value = 0.00594
$l1.Range("T3").Value = 0.00954
$l1.Range("U3").Value = 0.0216
$l1.Range("V3").Value = 0.0963

$l1.Range("M4").Value = 0.0213
$l1.Range("Q4").Value = 0.0213
$l1.Range("U4").Value = 0.0213

# ---------------------------------------------------------------------------
# Step 3 - refresh the percentages on the 5 division tabs
# ---------------------------------------------------------------------------

$ws = $wb.Worksheets.Item("IPS Clutches & Brakes Division")
$ws.Range("F2").Value = 0.0132
$ws.Range("I2").Value = 0.0286
$ws.Range("K2").Value = 0.0057
$ws.Range("L2").Value = 0.0105
$ws.Range("M2").Value = 0.0229
$ws.Range("Q2").Value = 0.0173
$ws.Range("V2").Value = 0.0866
$ws.Range("F3").Value = 0.01188
$ws.Range("I3").Value = 0.02574
$ws.Range("K3").Value = 0.00513
$ws.Range("L3").Value = 0.00945
$ws.Range("M3").Value = 0.02061
$ws.Range("Q3").Value = 0.01557
$ws.Range("V3").Value = 0.07794
$ws.Range("M4").Value = 0.020775
$ws.Range("Q4").Value = 0.020775
$ws.Range("U4").Value = 0.020775

$ws = $wb.Worksheets.Item("IPS Couplings Division")
$ws.Range("V2").Value = 0.1037
$ws.Range("V3").Value = 0.09333
$ws.Range("M4").Value = 0.02145
$ws.Range("Q4").Value = 0.02145
$ws.Range("U4").Value = 0.02145

$ws = $wb.Worksheets.Item("IPS Gearing Division")
$ws.Range("I2").Value = 0.0331
$ws.Range("L2").Value = 0.0173
$ws.Range("S2").Value = 0.0064
$ws.Range("T2").Value = 0.0161
$ws.Range("U2").Value = 0.0319
$ws.Range("V2").Value = 0.1406
$ws.Range("I3").Value = 0.02979
$ws.Range("L3").Value = 0.01557
$ws.Range("S3").Value = 0.00576
$ws.Range("T3").Value = 0.01449
$ws.Range("U3").Value = 0.02871
$ws.Range("V3").Value = 0.12654
$ws.Range("M4").Value = 0.0231
$ws.Range("Q4").Value = 0.0231
$ws.Range("U4").Value = 0.0231

$ws = $wb.Worksheets.Item("IPS Industrial Components Divi")
$ws.Range("N2").Value = 0.0252
$ws.Range("O2").Value = 0.0147
$ws.Range("Q2").Value = 0.0523
$ws.Range("U2").Value = 0.0454
$ws.Range("V2").Value = 0.1623
$ws.Range("N3").Value = 0.02268
$ws.Range("O3").Value = 0.01323
$ws.Range("Q3").Value = 0.04707
$ws.Range("U3").Value = 0.04086
$ws.Range("V3").Value = 0.14607
$ws.Range("M4").Value = 0.0162
$ws.Range("Q4").Value = 0.0162
$ws.Range("U4").Value = 0.0162

$ws = $wb.Worksheets.Item("IPS Segment Functions")
$ws.Range("F2").Value = 0.0102
$ws.Range("G2").Value = 0.0078
$ws.Range("H2").Value = 0.0094
$ws.Range("I2").Value = 0.0275
$ws.Range("J2").Value = 0.0078
$ws.Range("L2").Value = 0.0078
$ws.Range("M2").Value = 0.0249
$ws.Range("N2").Value = 0.007
$ws.Range("O2").Value = 0.007
$ws.Range("P2").Value = 0.0039
$ws.Range("Q2").Value = 0.0179
$ws.Range("R2").Value = 0.0062
$ws.Range("S2").Value = 0.007
$ws.Range("T2").Value = 0.0055
$ws.Range("U2").Value = 0.0187
$ws.Range("V2").Value = 0.0889
$ws.Range("F3").Value = 0.00918
$ws.Range("G3").Value = 0.00702
$ws.Range("H3").Value = 0.00846
$ws.Range("I3").Value = 0.02475
$ws.Range("J3").Value = 0.00702
$ws.Range("L3").Value = 0.00702
$ws.Range("M3").Value = 0.02241
$ws.Range("N3").Value = 0.0063
$ws.Range("O3").Value = 0.0063
$ws.Range("P3").Value = 0.00351
$ws.Range("Q3").Value = 0.01611
$ws.Range("R3").Value = 0.00558
$ws.Range("S3").Value = 0.0063
$ws.Range("T3").Value = 0.00495
$ws.Range("U3").Value = 0.01683
$ws.Range("V3").Value = 0.08001
$ws.Range("M4").Value = 0.022725
$ws.Range("Q4").Value = 0.022725
$ws.Range("U4").Value = 0.022725
